# CreateContact validData.xlsx update
# Replaces the 4 sample rows (rows 2-5) with 2 full contact rows (rows 1-2),
# refreshes hyperlinks, column widths, selection and page setup to match
# the new "Siddarth Sai" / Australia-UK contact data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to hold 4 partially-filled rows (rows 2-5); the new data
# only needs 2 fully-populated rows, so drop the extra rows first.
$ws.Rows("3:5").Delete()

# Remove the old hyperlink (it pointed at the old shyam@gmail.com in C2);
# new hyperlinks get (re)created below against the new data.
foreach ($hl in @($ws.Hyperlinks)) { $hl.Delete() }

# ---- Row 1 ----
$ws.Range("A1").Value = "Manoj"
$ws.Range("B1").Value = "Mummidi"
$ws.Range("C1").Value = "shyam@gmail.com"
$ws.Range("D1").Value = 1234567898
$ws.Range("E1").Value = "Road no:8A"
$ws.Range("F1").Value = "Sydney"
$ws.Range("G1").Value = "New South Wales"
$ws.Range("H1").Value = 500089
$ws.Range("I1").Value = "Australia"
$ws.Range("J1").Value = 2
$ws.Range("K1").Value = "March"
$ws.Range("L1").Value = 2020

# ---- Row 2 ----
$ws.Range("A2").Value = "Siddarth"
$ws.Range("B2").Value = "Sai"
$ws.Range("C2").Value = "saimanoj@yahoo.com"
$ws.Range("D2").Value = 9182820148
$ws.Range("E2").Value = "oxford street"
$ws.Range("F2").Value = "London"
$ws.Range("G2").Value = "North east"
$ws.Range("H2").Value = 530041
$ws.Range("I2").Value = "United Kingdom"
$ws.Range("J2").Value = 2
$ws.Range("K2").Value = "February"
$ws.Range("L2").Value = 1999

# ---- Hyperlinks (mailto:) on the email cells ----
# (apply the existing "Hyperlink" cell style by copying it from an already
# -styled cell, so no redundant style entries are added to styles.xml)
$ws.Hyperlinks.Add($ws.Range("C1"), "mailto:shyam@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:saimanoj@yahoo.com")
$hyperlinkStyle = $ws.Range("C2").Style
$ws.Range("C1").Style = $hyperlinkStyle
$ws.Range("C2").Style = $hyperlinkStyle

# ---- Column widths ----
$ws.Columns("C").ColumnWidth = 18.833333333333332
$ws.Columns("I").ColumnWidth = 15
$ws.Columns("L").ColumnWidth = 10

# ---- Selection / view ----
$ws.Range("H7").Select()

# ---- Page setup ----
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
